# The "Recorded By" column (G) contains values like "dnasr281@gmail.com, System"
# for sessions recorded by both the user and the System. This commit reorders
# those combined author lists so "System" is listed first:
#   "dnasr281@gmail.com, System"  ->  "System, dnasr281@gmail.com"
# Cells that only contain "dnasr281@gmail.com" (no "System") are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$colG = $ws.Range("G1:G" + $lastRow)
$colG.Replace("dnasr281@gmail.com, System", "System, dnasr281@gmail.com", -4163, 1, $false, $false, $false)
